$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.392.15'
$ws.Range('E2').Value = '  +0.54%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.525.27'
$ws.Range('E3').Value = '  +0.50%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '597.47'
$ws.Range('E5').Value = '  +0.49%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '173.80'
$ws.Range('E6').Value = '  +2.57%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  +3.20%  '
$ws.Range('E9').Value = '  +7.64%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.29'
$ws.Range('E10').Value = '  +0.14%  '
$ws.Range('E11').Value = '  -0.43%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.134.88'
$ws.Range('E12').Value = '  +0.54%  '
$ws.Range('E13').Value = '  +0.06%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.82'
$ws.Range('E14').Value = '  +1.92%  '
$ws.Range('E15').Value = '  +1.59%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '67.343.67'
$ws.Range('E16').Value = '  +0.52%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.562.81'
$ws.Range('E17').Value = '  +1.76%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.34'
$ws.Range('E18').Value = '  +0.27%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.32'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '397.74'
$ws.Range('E20').Value = '  +0.63%  '
$ws.Range('E21').Value = '  +0.61%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '73.50'
$ws.Range('E22').Value = '  +0.25%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.540'
$ws.Range('E23').Value = '  +1.68%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000123'
$ws.Range('E25').Value = '  -3.61%  '
$ws.Range('E27').Value = '  -0.61%  '
$ws.Range('E28').Value = '  -0.19%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.29'
$ws.Range('E29').Value = '  -1.21%  '
$ws.Range('E30').Value = '  -0.97%  '
$ws.Range('E31').Value = '  +0.62%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '24.17'
$ws.Range('E32').Value = '  +2.43%  '
$ws.Range('E33').Value = '  -0.99%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.65'
$ws.Range('E34').Value = '  +2.80%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '163.30'
$ws.Range('E35').Value = '  +0.40%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.898'
$ws.Range('E36').Value = '  -0.70%  '
$ws.Range('E37').Value = '  -1.21%  '
$ws.Range('E38').Value = '  +3.92%  '
$ws.Range('E39').Value = '  +1.24%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '27.70'
$ws.Range('E40').Value = '  +4.13%  '
$ws.Range('E41').Value = '  -0.95%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '26.50'
$ws.Range('E42').Value = '  +0.32%  '
$ws.Range('E43').Value = '  +2.54%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.802.40'
$ws.Range('E44').Value = '  -1.53%  '
$ws.Range('E45').Value = '  -1.34%  '
$ws.Range('E46').Value = '  -2.63%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '341.88'
$ws.Range('E47').Value = '  -2.44%  '
$ws.Range('E48').Value = '  +1.30%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '33.80'
$ws.Range('E49').Value = '  +0.26%  '
$ws.Range('E50').Value = '  +0.05%  '
$ws.Range('E51').Value = '  -0.74%  '
